$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2436
$ws.Range("I33").Value = 2771
$ws.Range("J33").Value = 91
$ws.Range("K33").Value = 2771
$ws.Range("L33").Value = 91
$ws.Range("M33").Value = -2542
$ws.Range("N33").Value = -549
$ws.Range("H135").Value = 3914.138
$ws.Range("I135").Value = 2729.6667
$ws.Range("K135").Value = 24567.0003
$ws.Range("M135").Value = -22032.0003
$ws.Range("H137").Value = 2623.543
$ws.Range("I137").Value = 1847.52
$ws.Range("J137").Value = 4563.6
$ws.Range("K137").Value = 5542.559999999999
$ws.Range("L137").Value = 13690.8
$ws.Range("M137").Value = -2992.559999999999
$ws.Range("N137").Value = -18790.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27350.916
$ws.Range("I32").Value = 37356.875
$ws.Range("K32").Value = 37356.875
$ws.Range("M32").Value = -37069.875
$ws.Range("H61").Value = 3945.6492
$ws.Range("I61").Value = 2250.6924
$ws.Range("K61").Value = 2250.6924
$ws.Range("M61").Value = -2038.6924
$ws.Range("H124").Value = 18414.5
$ws.Range("J124").Value = 18414.5
$ws.Range("L124").Value = 18414.5
$ws.Range("N124").Value = -28234.5
$ws.Range("H132").Value = 4399.3335
$ws.Range("I132").Value = 4103
$ws.Range("J132").Value = 5288.3335
$ws.Range("K132").Value = 12309
$ws.Range("L132").Value = 15865.0005
$ws.Range("M132").Value = -9779
$ws.Range("N132").Value = -20925.0005
$ws.Range("H136").Value = 3945.6492
$ws.Range("I136").Value = 2250.6924
$ws.Range("K136").Value = 6752.0772
$ws.Range("M136").Value = -4202.0772

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 54113.824
$ws.Range("I86").Value = 1335.75
$ws.Range("J86").Value = 101027.664
$ws.Range("K86").Value = 1335.75
$ws.Range("L86").Value = 101027.664
$ws.Range("M86").Value = -212.75
$ws.Range("N86").Value = -103273.664
$ws.Range("H89").Value = 54113.824
$ws.Range("I89").Value = 1335.75
$ws.Range("J89").Value = 101027.664
$ws.Range("K89").Value = 6678.75
$ws.Range("L89").Value = 505138.32
$ws.Range("M89").Value = -1062.75
$ws.Range("N89").Value = -516370.32
$ws.Range("H134").Value = 3719.5
$ws.Range("I134").Value = 2989.5715
$ws.Range("J134").Value = 6785.2
$ws.Range("K134").Value = 8968.7145
$ws.Range("L134").Value = 20355.6
$ws.Range("M134").Value = -6433.7145
$ws.Range("N134").Value = -25425.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 6000
$ws.Range("I39").Value = 6000
$ws.Range("K39").Value = 6000
$ws.Range("M39").Value = -5609
$ws.Range("H49").Value = 6000
$ws.Range("I49").Value = 6000
$ws.Range("K49").Value = 6000
$ws.Range("M49").Value = -5818
$ws.Range("H58").Value = 13945.333
$ws.Range("I58").Value = 1921
$ws.Range("J58").Value = 28975.75
$ws.Range("K58").Value = 1921
$ws.Range("L58").Value = 28975.75
$ws.Range("M58").Value = -1718
$ws.Range("N58").Value = -29381.75
$ws.Range("H132").Value = 15566
$ws.Range("I132").Value = 8129.4
$ws.Range("J132").Value = 23002.6
$ws.Range("K132").Value = 24388.2
$ws.Range("L132").Value = 69007.79999999999
$ws.Range("M132").Value = -21858.2
$ws.Range("N132").Value = -74067.79999999999
$ws.Range("H136").Value = 13945.333
$ws.Range("I136").Value = 1921
$ws.Range("J136").Value = 28975.75
$ws.Range("K136").Value = 5763
$ws.Range("L136").Value = 86927.25
$ws.Range("M136").Value = -3213
$ws.Range("N136").Value = -92027.25
$ws.Range("H141").Value = 148542.7
$ws.Range("J141").Value = 148542.7
$ws.Range("L141").Value = 148542.7
$ws.Range("N141").Value = -158902.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 493.73914
$ws.Range("I5").Value = 468.14285
$ws.Range("J5").Value = 533.55554
$ws.Range("K5").Value = 1404.42855
$ws.Range("L5").Value = 1600.66662
$ws.Range("M5").Value = -1292.42855
$ws.Range("N5").Value = -1824.66662
$ws.Range("H133").Value = 7426
$ws.Range("I133").Value = 7426
$ws.Range("K133").Value = 22278
$ws.Range("M133").Value = -17218
$ws.Range("H135").Value = 493.73914
$ws.Range("I135").Value = 468.14285
$ws.Range("J135").Value = 533.55554
$ws.Range("K135").Value = 4213.28565
$ws.Range("L135").Value = 4801.99986
$ws.Range("M135").Value = -1678.28565
$ws.Range("N135").Value = -9871.99986
$ws.Range("H136").Value = 4454
$ws.Range("I136").Value = 4298
$ws.Range("K136").Value = 12894
$ws.Range("M136").Value = -7794
$ws.Range("H139").Value = 3848
$ws.Range("I139").Value = 4644
$ws.Range("K139").Value = 13932
$ws.Range("M139").Value = -8792
$ws.Range("H141").Value = 2728.1875
$ws.Range("I141").Value = 2510.0667
$ws.Range("J141").Value = 6000
$ws.Range("K141").Value = 7530.2001
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = -2350.2001
$ws.Range("N141").Value = -28360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7087.1177
$ws.Range("I80").Value = 4525
$ws.Range("K80").Value = 4525
$ws.Range("M80").Value = -3527
$ws.Range("H83").Value = 7087.1177
$ws.Range("I83").Value = 4525
$ws.Range("K83").Value = 22625
$ws.Range("M83").Value = -17633
$ws.Range("H109").Value = 25070.072
$ws.Range("J109").Value = 25070.072
$ws.Range("L109").Value = 25070.072
$ws.Range("N109").Value = -27150.072
$ws.Range("H132").Value = 3083254
$ws.Range("I132").Value = 3339400.2
$ws.Range("K132").Value = 10018200.6
$ws.Range("M132").Value = -10015670.6

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1934.1428
$ws.Range("I55").Value = 570.0714
$ws.Range("K55").Value = 570.0714
$ws.Range("M55").Value = -397.0714
$ws.Range("H61").Value = 3010.7896
$ws.Range("I61").Value = 3010.7896
$ws.Range("K61").Value = 3010.7896
$ws.Range("M61").Value = -2808.7896
$ws.Range("H68").Value = 5041.6665
$ws.Range("I68").Value = 2819
$ws.Range("J68").Value = 7820
$ws.Range("K68").Value = 2819
$ws.Range("L68").Value = 7820
$ws.Range("M68").Value = -2070
$ws.Range("N68").Value = -9318
$ws.Range("H71").Value = 5041.6665
$ws.Range("I71").Value = 2819
$ws.Range("J71").Value = 7820
$ws.Range("K71").Value = 14095
$ws.Range("L71").Value = 39100
$ws.Range("M71").Value = -10351
$ws.Range("N71").Value = -46588
$ws.Range("H113").Value = 3010.7896
$ws.Range("I113").Value = 3010.7896
$ws.Range("K113").Value = 3010.7896
$ws.Range("M113").Value = -840.7896000000001
$ws.Range("H132").Value = 3657.9678
$ws.Range("J132").Value = 6131.5
$ws.Range("L132").Value = 18394.5
$ws.Range("N132").Value = -23454.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7295.75
$ws.Range("J132").Value = 13963.333
$ws.Range("L132").Value = 41889.999
$ws.Range("M132").Value = -16499.7149
$ws.Range("N132").Value = -46949.999
$ws.Range("H136").Value = 4513.294
$ws.Range("I136").Value = 4609.4287
$ws.Range("K136").Value = 13828.2861
$ws.Range("M136").Value = -11278.2861
